# "changed to new battery": update the Battery Specs block on Sheet1
# (column I) to reflect the new battery's capacity and max draw.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "5000mAh"
$ws.Range("I4").Value = "500A Max Draw"

# Leave the cursor where the author left it after making the edit.
$ws.Range("J7").Select()
